# Fruta / hortaliza, semanal
# Insert a new weekly record at row 171, shifting subsequent rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 171 (existing rows 171..219 shift to 172..220)
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row with the new record's data
$ws.Range("A171").Value = 7
$ws.Range("B171").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C171").Value = "Ñuble"
$ws.Range("D171").Value = 44663
$ws.Range("E171").Value = 16
$ws.Range("F171").Value = 100112043
$ws.Range("G171").Value = "Pepino ensalada"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 120
$ws.Range("K171").Value = 16000
$ws.Range("L171").Value = 16500
$ws.Range("M171").Value = 16250
$ws.Range("N171").Value = "$/caja 80 unidades"
$ws.Range("O171").Value = "Región del Maule"
$ws.Range("P171").Value = 203
$ws.Range("Q171").Value = 80
$ws.Range("R171").Value = "Hortaliza"
